$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1.8315
$ws.Range("G2").Value = -0.5327470492941448
$ws.Range("H2").Value = -1.006711409395973
$ws.Range("I2").Value = -1.02985420041657
$ws.Range("J2").Value = -1.02985420041657
$ws.Range("K2").Value = -17.5
$ws.Range("L2").Value = -1.349996142868163
$ws.Range("U2").Value = 4.25
$ws.Range("V2").Value = 0.0154320987654321
$ws.Range("W2").Value = -1.174894217207334
$ws.Range("X2").Value = 0.07704211110163209
$ws.Range("Y2").Value = -1.251936328308966
$ws.Range("Z2").Value = 0.4917305212047645
$ws.Range("AA2").Value = -1.379435191889934
$ws.Range("AB2").Value = 0.07689136991668613
$ws.Range("AC2").Value = -1.465267093811619
$ws.Range("AD2").Value = 7.359
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 7.359
$ws.Range("AG2").Value = 3.109
$ws.Range("AH2").Value = 0.02602569679479699
$ws.Range("AI2").Value = 0.2525481313703284
$ws.Range("AJ2").Value = 0.01116301448068106
$ws.Range("AK2").Value = 0.1249146209168709
$ws.Range("AL2").Value = 0.872
$ws.Range("AM2").Value = 0.583
$ws.Range("AN2").Value = -0.5639080459770115
$ws.Range("AO2").Value = -15.30963302752293
$ws.Range("AP2").Value = -0.2382375478927203
$ws.Range("AQ2").Value = -22.89879931389365

# Row 3
$ws.Range("D3").Value = -0.148
$ws.Range("G3").Value = -2.078947368421053
$ws.Range("H3").Value = -5.008771929824562
$ws.Range("I3").Value = -5.052631578947369
$ws.Range("J3").Value = -5.052631578947369
$ws.Range("K3").Value = -5.39
$ws.Range("L3").Value = -4.728070175438597
$ws.Range("U3").Value = 1.12
$ws.Range("V3").Value = 0.009443507588532884
$ws.Range("W3").Value = -0.2130434782608696
$ws.Range("X3").Value = 0.07690773952214357
$ws.Range("Y3").Value = -0.2899512177830131
$ws.Range("Z3").Value = 0.04863481228668941
$ws.Range("AA3").Value = -0.2457337883959044
$ws.Range("AB3").Value = 0.07689136991668613
$ws.Range("AC3").Value = -0.3226251583125905
$ws.Range("AD3").Value = 0.033
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.033
$ws.Range("AG3").Value = -1.087
$ws.Range("AH3").Value = 0.0002781688063186466
$ws.Range("AI3").Value = 0.001733830715073819
$ws.Range("AJ3").Value = -0.009250040421059799
$ws.Range("AK3").Value = -0.06068218612181098
$ws.Range("AL3").Value = 0.007
$ws.Range("AM3").Value = -0.282
$ws.Range("AN3").Value = -0.005779334500875657
$ws.Range("AO3").Value = -822.8571428571428
$ws.Range("AP3").Value = 0.1903677758318739
$ws.Range("AQ3").Value = 20.42553191489362

# Row 4
$ws.Range("G4").Value = -0.2920338983050847
$ws.Range("H4").Value = -0.3084745762711864
$ws.Range("I4").Value = -0.3228813559322034
$ws.Range("J4").Value = -0.3228813559322034
$ws.Range("K4").Value = -8.33
$ws.Range("L4").Value = -0.7059322033898304
$ws.Range("U4").Value = 1.4
$ws.Range("V4").Value = 0.04697986577181208
$ws.Range("W4").Value = -1.174894217207334
$ws.Range("X4").Value = 0.0878979725699162
$ws.Range("Y4").Value = -1.26279218977725
$ws.Range("Z4").Value = 4.272266473569875
$ws.Range("AA4").Value = -1.379435191889934
$ws.Range("AB4").Value = 0.08583190192168536
$ws.Range("AC4").Value = -1.465267093811619
$ws.Range("AD4").Value = 6.93
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 6.93
$ws.Range("AG4").Value = 5.529999999999999
$ws.Range("AH4").Value = 0.1886741083582902
$ws.Range("AI4").Value = 0.8067520372526193
$ws.Range("AJ4").Value = 0.1565242003962638
$ws.Range("AK4").Value = 0.7691237830319888
$ws.Range("AL4").Value = 0.855
$ws.Range("AM4").Value = 0.855
$ws.Range("AN4").Value = -1.903846153846154
$ws.Range("AO4").Value = -4.456140350877194
$ws.Range("AP4").Value = -1.519230769230769
$ws.Range("AQ4").Value = -4.456140350877194

# Row 5
$ws.Range("D5").Value = 3.811
$ws.Range("G5").Value = -47.3913043478261
$ws.Range("H5").Value = -160.8695652173913
$ws.Range("I5").Value = -164.3478260869565
$ws.Range("J5").Value = -164.3478260869565
$ws.Range("K5").Value = -3.78
$ws.Range("L5").Value = -164.3478260869565
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 0.01362204724409449
$ws.Range("W5").Value = -2.25
$ws.Range("X5").Value = 0.07704211110163209
$ws.Range("Y5").Value = -2.327042111101632
$ws.Range("Z5").Value = 0.1437500000000001
$ws.Range("AA5").Value = -23.62500000000001
$ws.Range("AB5").Value = 0.0768587700548469
$ws.Range("AC5").Value = -23.70185877005486
$ws.Range("AD5").Value = 0.396
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0.396
$ws.Range("AG5").Value = -1.334
$ws.Range("AH5").Value = 0.003108417846714183
$ws.Range("AI5").Value = 0.2612137203166227
$ws.Range("AJ5").Value = -0.01061544093072112
$ws.Range("AK5").Value = 6.233644859813086
$ws.Range("AL5").Value = 0.01
$ws.Range("AM5").Value = 0.01
$ws.Range("AN5").Value = -0.107027027027027
$ws.Range("AO5").Value = -378
$ws.Range("AP5").Value = 0.3605405405405405
$ws.Range("AQ5").Value = -378
